# Apply the updated "想去人数" (interested-count) values to the
# "展览" (sheet 1) and "全部类型" (sheet 4) worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 56
$ws1.Range("F8").Value  = 81
$ws1.Range("F9").Value  = 8528
$ws1.Range("F13").Value = 925
$ws1.Range("F14").Value = 87
$ws1.Range("F16").Value = 225
$ws1.Range("F17").Value = 201
$ws1.Range("F20").Value = 968

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 56
$ws4.Range("F10").Value = 81
$ws4.Range("F11").Value = 8528
$ws4.Range("F15").Value = 925
$ws4.Range("F16").Value = 87
$ws4.Range("F18").Value = 225
$ws4.Range("F19").Value = 201
$ws4.Range("F22").Value = 968
